$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "bleu" = "noir"
    "pas de résultat ni de publication" = "pas de résultat postés ni publiés"
    "résultat et / ou publication posté dans les 36 mois" = "résultat postés ou publiés dans les 36 mois"
    "résultat et / ou publication posté dans les 12 mois" = "résultat postés ou publiés dans les 12 mois"
    "résultat et / ou publication posté" = "résultat postés ou publiés"
}

$used = $ws.UsedRange
$nrows = $used.Rows.Count
$ncols = $used.Columns.Count

for ($r = 1; $r -le $nrows; $r++) {
    for ($c = 1; $c -le $ncols; $c++) {
        $cell = $used.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($null -ne $val -and $replacements.ContainsKey($val)) {
            $cell.Value2 = $replacements[$val]
        }
    }
}
